{"js": "// Add a \"Conclusion:\" heading + summary paragraph, right before the final\n// (empty) paragraph that carries the _GoBack bookmark at the end of the body.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// Insert the summary body paragraph first, directly before the trailing\n// empty paragraph. It inherits that paragraph's (default/Normal) style, so\n// no explicit style override is required.\nconst bodyParagraph = lastParagraph.insertParagraph(\n  \"GitHub is a powerful tool for version control and team collaboration. It helps manage projects efficiently by tracking every change made by each member.\",\n  \"Before\"\n);\n\n// Insert the \"Conclusion:\" heading right before that new summary paragraph.\nconst headingParagraph = bodyParagraph.insertParagraph(\"Conclusion:\", \"Before\");\nheadingParagraph.styleBuiltIn = Word.Style.heading2;\n\nawait context.sync();\n", "ps1": "# Add a \"Conclusion:\" heading + summary paragraph, right before the final\n# (empty) paragraph that carries the _GoBack bookmark at the end of the body.\n$d = $word.ActiveDocument\n\n$lastParagraph = $d.Paragraphs.Last\n\n# Insert the summary body paragraph first, directly before the trailing\n# empty paragraph. It naturally picks up the default/Normal style, so no\n# explicit style assignment is required.\n$lastParagraph.Range.InsertParagraphBefore()\n$bodyParagraph = $d.Paragraphs.Item($d.Paragraphs.Count - 1)\n$bodyParagraph.Range.Text = \"GitHub is a powerful tool for version control and team collaboration. It helps manage projects efficiently by tracking every change made by each member.\"\n\n# Insert the \"Conclusion:\" heading right before that new summary paragraph.\n$bodyParagraph.Range.InsertParagraphBefore()\n$headingParagraph = $d.Paragraphs.Item($d.Paragraphs.Count - 2)\n$headingParagraph.Range.Text = \"Conclusion:\"\n$headingParagraph.Style = \"Heading 2\"\n"}
